# This script removes two blocks of paragraphs from chapter 3.1.3:
#
# 1) Right after the "Output" heading, the paragraph
#    "Uttrekket inneholder ETT arkiv og EN arkivdeler og er riktig
#    avsluttet.", the following empty paragraph, and the following
#    "AND/OR" paragraph are removed (the empty paragraph that used to
#    sit right after "AND/OR" is kept, so "Output" is now followed by a
#    single empty paragraph and then the "ANTALL arkivdeler" text).
#
# 2) After the table, the "AND/OR" paragraph, the following empty
#    paragraph, and the "Arkivdelstatus er satt til ..." paragraph are
#    removed (the empty paragraph that follows the table, and the one
#    that used to follow the removed "Arkivdelstatus ..." paragraph,
#    are both kept).
#
# We locate the anchor paragraphs by their text and delete the ranges
# that span from the start of the first paragraph to be removed through
# the end of the last one in each block. We process the later block
# first so earlier paragraph indices/ranges stay valid.

$d = $word.ActiveDocument

function Find-ParagraphIndex($doc, $text, $exact) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text
        if ($exact) {
            if ($t -eq $text) { return $i }
        } else {
            if ($t -like "*$text*") { return $i }
        }
    }
    return -1
}

# --- Block 2 (later in the document): remove "AND/OR" (after table),
#     the blank paragraph after it, and the "Arkivdelstatus er satt
#     til" paragraph. Process this one first (it's further down in the
#     document) so the earlier block's indices aren't disturbed.

$statusIdx = Find-ParagraphIndex $d "Arkivdelstatus er satt til" $false
if ($statusIdx -eq -1) {
    throw "Could not find 'Arkivdelstatus er satt til' paragraph"
}

# The "AND/OR" paragraph immediately preceding the blank paragraph that
# precedes the status paragraph is two paragraphs before it.
$andOr2Idx = $statusIdx - 2
$andOr2Para = $d.Paragraphs.Item($andOr2Idx)
if ($andOr2Para.Range.Text -notlike "AND/OR*") {
    throw "Unexpected paragraph content before status paragraph: $($andOr2Para.Range.Text)"
}

$start2 = $andOr2Para.Range.Start
$end2 = $d.Paragraphs.Item($statusIdx).Range.End
$d.Range($start2, $end2).Delete()

# --- Block 1: remove "Uttrekket inneholder ETT arkiv ..." paragraph,
#     the blank paragraph after it, and the following "AND/OR"
#     paragraph.

$uttrekketIdx = Find-ParagraphIndex $d "Uttrekket inneholder" $false
if ($uttrekketIdx -eq -1) {
    throw "Could not find 'Uttrekket inneholder' paragraph"
}
$uttrekketPara = $d.Paragraphs.Item($uttrekketIdx)
if ($uttrekketPara.Range.Text -notlike "*ETT arkiv*") {
    throw "Unexpected first 'Uttrekket inneholder' paragraph: $($uttrekketPara.Range.Text)"
}

# The matching "AND/OR" paragraph is two paragraphs after it.
$andOr1Idx = $uttrekketIdx + 2
$andOr1Para = $d.Paragraphs.Item($andOr1Idx)
if ($andOr1Para.Range.Text -notlike "AND/OR*") {
    throw "Unexpected paragraph after 'Uttrekket inneholder' block: $($andOr1Para.Range.Text)"
}

$start1 = $uttrekketPara.Range.Start
$end1 = $andOr1Para.Range.End
$d.Range($start1, $end1).Delete()
